{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Change described by the diff:\n//   1. Delete the table row whose first cell reads \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\"\n//      (second cell held \"5\").\n//   2. Append a brand-new row at the end of the (same) table whose first\n//      cell reads \"new row name\" and whose second cell is left empty.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in the document body.\");\n}\n\nconst table = tables.items[0];\n\n// --- 1. Locate and delete the \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\" row -------------------\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nfor (const row of rows) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nconst targetLabel = \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\";\nlet rowToDelete = null;\nfor (const row of rows) {\n  const firstCellText = row.cells.items[0] && row.cells.items[0].body.text\n    ? row.cells.items[0].body.text.trim()\n    : \"\";\n  if (firstCellText === targetLabel) {\n    rowToDelete = row;\n    break;\n  }\n}\n\nif (rowToDelete) {\n  rowToDelete.delete();\n  await context.sync();\n}\n\n// --- 2. Append the new row at the end of the table --------------------\ntable.addRows(\"End\", 1, [[\"new row name\", \"\"]]);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Change described by the diff:\n#   1. Delete the table row whose first cell reads \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\"\n#      (second cell held \"5\").\n#   2. Append a brand-new row at the end of the (same) table whose first\n#      cell reads \"new row name\" and whose second cell is left empty.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- 1. Locate and delete the \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\" row --------------------\n$targetLabel = \"\u041a\u043e\u043b\u0438\u0447\u0435\u0441\u0442\u0432\u043e \u043a\u0440\u0435\u0434\u0438\u0442\u043e\u0432\"\n$targetIndex = -1\n$rowCount = $t.Rows.Count\nfor ($i = 1; $i -le $rowCount; $i++) {\n    $row = $t.Rows.Item($i)\n    # Cell text carries a trailing cell-mark (chr 7) / paragraph mark (chr 13)\n    $cellText = $row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13)\n    if ($cellText -eq $targetLabel) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $t.Rows.Item($targetIndex).Delete()\n}\n\n# --- 2. Append the new row at the end of the table --------------------\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"new row name\"\n"}
